$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header formatting tweaks (match "insee"/"num_siret" text style, and
#     give "hauteur_max" its own text style) ---
$ws.Range("C2").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R2").Copy()
$ws.Range("R1").PasteSpecial(-4122)   # xlPasteFormats

# hauteur_max column becomes a text column (id 75114-P-001 currently stores "5")
$ws.Range("Q1").NumberFormat = "@"
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "290"

# --- Duplicate row 2 into a new row 3 (second parking record) ---
$ws.Range("A2:AD2").Copy($ws.Range("A3:AD3"))

# New record: same parking (REPUBLIQUE) but a different id / ouvrage
$ws.Range("A3").Value = "75114-P-002"

# Non applicable height -> "N/A" (was "NA") for this second record
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "N/A"

# --- Sheet-level touches to mirror the authored workbook ---
$ws.Range("A1:AD2").AutoFilter() | Out-Null
$ws.Range("R3").Select() | Out-Null
